$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.328.07"
$ws.Range("E2").Value = "'  +0.50%  "
$ws.Range("D3").Value = "'1.873.18"
$ws.Range("E3").Value = "'  +0.60%  "
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("D5").Value = "'0.7126"
$ws.Range("E5").Value = "'  +0.71%  "
$ws.Range("D6").Value = "'241.50"
$ws.Range("E6").Value = "'  +0.25%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("B8").Value = "'Dogecoin"
$ws.Range("C8").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07795"
$ws.Range("E8").Value = "'  +1.93%  "
$ws.Range("B9").Value = "'Cardano"
$ws.Range("C9").Value = "'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.3112"
$ws.Range("E9").Value = "'  +0.86%  "
$ws.Range("D10").Value = "'25.11"
$ws.Range("E10").Value = "'  +1.80%  "
$ws.Range("E11").Value = "'  +1.01%  "
$ws.Range("D12").Value = "'1.863.23"
$ws.Range("E12").Value = "'  +0.64%  "
$ws.Range("D13").Value = "'5.236"
$ws.Range("D14").Value = "'0.7120"
$ws.Range("E14").Value = "'  +0.34%  "
$ws.Range("D15").Value = "'91.12"
$ws.Range("E15").Value = "'  -0.10%  "
$ws.Range("D16").Value = "'29.337.38"
$ws.Range("E16").Value = "'  +0.44%  "
$ws.Range("D17").Value = "'6.093"
$ws.Range("E17").Value = "'  +2.97%  "
$ws.Range("D18").Value = "'0.000008231"
$ws.Range("E18").Value = "'  +5.33%  "
$ws.Range("D19").Value = "'240.05"
$ws.Range("E19").Value = "'  -1.07%  "
$ws.Range("D20").Value = "'13.21"
$ws.Range("E20").Value = "'  +1.04%  "
$ws.Range("D21").Value = "'2.122.42"
$ws.Range("E21").Value = "'  +0.36%  "
$ws.Range("E22").Value = "'  -0.02%  "
$ws.Range("D23").Value = "'7.756"
$ws.Range("E23").Value = "'  -1.27%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "'  -0.04%  "
$ws.Range("D25").Value = "'0.1599"
$ws.Range("E25").Value = "'  +0.44%  "
$ws.Range("D26").Value = "'162.89"
$ws.Range("E26").Value = "'  -0.32%  "
$ws.Range("D27").Value = "'9.035"
$ws.Range("E27").Value = "'  +1.12%  "
$ws.Range("E28").Value = "'  +0.12%  "
$ws.Range("E29").Value = "'  +0.78%  "
$ws.Range("D30").Value = "'4.417"
$ws.Range("E30").Value = "'  +0.27%  "
$ws.Range("E31").Value = "'  -1.85%  "
$ws.Range("D32").Value = "'4.316"
$ws.Range("E32").Value = "'  +2.27%  "
$ws.Range("D33").Value = "'0.05295"
$ws.Range("E33").Value = "'  +3.22%  "
$ws.Range("E34").Value = "'  +1.48%  "
$ws.Range("D35").Value = "'1.177"
$ws.Range("E35").Value = "'  +1.38%  "
$ws.Range("D36").Value = "'0.7429"
$ws.Range("E36").Value = "'  -6.82%  "
$ws.Range("D37").Value = "'2.700"
$ws.Range("E37").Value = "'  +0.78%  "
$ws.Range("D38").Value = "'0.01872"
$ws.Range("E38").Value = "'  +1.75%  "
$ws.Range("D39").Value = "'1.220.14"
$ws.Range("E39").Value = "'  +4.58%  "
$ws.Range("D40").Value = "'2.728"
$ws.Range("E40").Value = "'  +1.14%  "
$ws.Range("D41").Value = "'6.553"
$ws.Range("E41").Value = "'  +6.06%  "
$ws.Range("D42").Value = "'110.86"
$ws.Range("E42").Value = "'  +8.54%  "
$ws.Range("D43").Value = "'0.8867"
$ws.Range("E43").Value = "'  -0.19%  "
$ws.Range("D44").Value = "'72.67"
$ws.Range("E44").Value = "'  -0.21%  "
$ws.Range("D45").Value = "'1.0000"
$ws.Range("E45").Value = "'  -0.01%  "
$ws.Range("D46").Value = "'2.019.04"
$ws.Range("E46").Value = "'  +0.29%  "
$ws.Range("D47").Value = "'1.802"
$ws.Range("E47").Value = "'  +1.97%  "
$ws.Range("D48").Value = "'0.5185"
$ws.Range("E48").Value = "'  +0.03%  "
$ws.Range("D50").Value = "'9.400"
$ws.Range("E50").Value = "'  +1.01%  "
$ws.Range("E51").Value = "'  +1.14%  "
